# Final project-data edits to the Gantt chart ("final final commit! PROJECT DONE")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# --- Row 10: progress bumped to 100% ---
$ws.Range("D10").Value = 1

# --- Row 11: task renamed ---
$ws.Range("B11").Value = "DTO Object Controller"

# --- Row 13: progress bumped to 100%, actual hours recorded ---
$ws.Range("D13").Value = 1
$ws.Range("G13").Value = "2 hours"
$ws.Range("H13").Value = "4 hours"

# --- Row 14: sample placeholder task cleared out ---
$ws.Range("B14").Value = ""

# --- Row 16: new task ---
$ws.Range("B16").Value = "Place design into box"
$ws.Range("C16").Value = "Matt"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 44465
$ws.Range("F16").Value = 44469
$ws.Range("G16").Value = "1 hours"
$ws.Range("H16").Value = "2 hours"

# --- Row 17: new task ---
$ws.Range("B17").Value = "Backend RaspPI program"
$ws.Range("C17").Value = "Henry"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 44458
$ws.Range("F17").Value = 44468
$ws.Range("G17").Value = "48 hours"
$ws.Range("H17").Value = "60 hours"

# --- Row 18: new task ---
$ws.Range("B18").Value = "Database setup"
$ws.Range("C18").Value = "Brendan/Matt"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 44467
$ws.Range("F18").Value = 44469
$ws.Range("G18").Value = "1 hour"
$ws.Range("H18").Value = "5.5 hours"

# --- Row 19: new task (was sample "Task 4") ---
$ws.Range("B19").Value = "Solder to 2m cable/ connections"
$ws.Range("C19").Value = "Matt"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 44469
$ws.Range("F19").Value = 44469
$ws.Range("G19").Value = "2 hours"
$ws.Range("H19").Value = "1.5 hours"

# --- Row 20: new task (was sample "Task 5") ---
$ws.Range("B20").Value = "Finish front end"
$ws.Range("C20").Value = "Brenden/Brendan"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 44467
$ws.Range("F20").Value = 44469
$ws.Range("G20").Value = "20 hours"
$ws.Range("H20").Value = "15 hours"

# --- View state: selection moved to B12 ---
$ws.Range("B12").Select()

# --- Column width tweaks (task name / assigned-to columns) ---
$ws.Columns("B").ColumnWidth = 30.88671875
$ws.Columns("C").ColumnWidth = 23.5546875
